# Correct the "SmithInc" company name (missing space) to "Smith Inc"
# for the John Smith row in the customer sample sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Select()
$ws.Range("C2").Value = "Smith Inc"
